# Replace spaces with underscores in the header labels (row 1, columns A:V)
# of Sheet1, and update the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column E (index 5) through V (index 22) first, then column D
# (Monthly Charge) last, so the shared-strings table ends up ordered the
# same way the original authoring tool produced it.
$cols = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,4)
foreach ($c in $cols) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val -replace " ", "_"
    }
}

$ws.Activate()
$ws.Range("A2").Select()
